# Update cryptos list data (prices / 1h volume change) and swap
# the Stacks / NEARProtocol rows, per the Feb 14 2024 GitHub Actions run.
#
# Note: price values in column D are stored as plain text (e.g.
# "49.583.59", "1.00", "0.547") rather than numbers. Setting .Value
# directly with a numeric-looking string would make Excel auto-convert
# it to a real number (losing formatting, e.g. "1.00" -> 1). To avoid
# that, the cell's number format is temporarily switched to Text ("@")
# before the value is written, then the style is reset back to Normal
# so the resulting cell keeps the original (default/General) styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Row 2: Bitcoin ---
Set-TextValue "D2" "49.583.59"
$ws.Range("E2").Value = "  -0.55%  "

# --- Row 3: Ethereum ---
Set-TextValue "D3" "2.641.03"
$ws.Range("E3").Value = "  -0.26%  "

# --- Row 4: TetherUSD ---
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.00%  "

# --- Row 5: Solana ---
Set-TextValue "D5" "111.90"
$ws.Range("E5").Value = "  -1.67%  "

# --- Row 6: BNB ---
Set-TextValue "D6" "325.76"
$ws.Range("E6").Value = "  -0.41%  "

# --- Row 7: XRP ---
$ws.Range("E7").Value = "  -1.03%  "

# --- Row 8: USDC ---
$ws.Range("E8").Value = "  +0.00%  "

# --- Row 9: Cardano ---
Set-TextValue "D9" "0.547"
$ws.Range("E9").Value = "  -1.29%  "

# --- Row 10: Avalanche ---
Set-TextValue "D10" "39.57"
$ws.Range("E10").Value = "  -3.77%  "

# --- Row 11: Chainlink ---
$ws.Range("E11").Value = "  -0.39%  "

# --- Row 12: Dogecoin ---
$ws.Range("E12").Value = "  -1.41%  "

# --- Row 13: TRON (unchanged) ---

# --- Row 14: Polkadot ---
Set-TextValue "D14" "7.52"
$ws.Range("E14").Value = "  +2.39%  "

# --- Row 15: WrappedliquidstakedEther2.0 ---
Set-TextValue "D15" "3.053.14"
$ws.Range("E15").Value = "  -0.38%  "

# --- Row 16: WrappedEther ---
Set-TextValue "D16" "2.636.17"
$ws.Range("E16").Value = "  -0.27%  "

# --- Row 17: Polygon ---
$ws.Range("E17").Value = "  -2.15%  "

# --- Row 18: WrappedBTC ---
Set-TextValue "D18" "49.537.56"
$ws.Range("E18").Value = "  -0.54%  "

# --- Row 19: InternetComputer(DFINITY) ---
Set-TextValue "D19" "13.11"
$ws.Range("E19").Value = "  -0.48%  "

# --- Row 20: ImmutableX ---
$ws.Range("E20").Value = "  -0.58%  "

# --- Row 21: Uniswap ---
Set-TextValue "D21" "6.68"
$ws.Range("E21").Value = "  -1.15%  "

# --- Row 22: ShibaInu ---
$ws.Range("E22").Value = "  -0.82%  "

# --- Row 23: BitcoinCash ---
Set-TextValue "D23" "268.72"
$ws.Range("E23").Value = "  -2.87%  "

# --- Row 24: Litecoin ---
Set-TextValue "D24" "69.07"
$ws.Range("E24").Value = "  -4.31%  "

# --- Row 25: PancakeSwap ---
Set-TextValue "D25" "2.56"
$ws.Range("E25").Value = "  -1.17%  "

# --- Row 26: EthereumClassic ---
Set-TextValue "D26" "26.08"
$ws.Range("E26").Value = "  -2.53%  "

# --- Row 27: Dai ---
$ws.Range("E27").Value = "  +0.05%  "

# --- Row 28: Cosmos ---
Set-TextValue "D28" "10.20"
$ws.Range("E28").Value = "  +1.66%  "

# --- Row 29: Toncoin ---
Set-TextValue "D29" "2.20"
$ws.Range("E29").Value = "  +0.19%  "

# --- Row 30: Kaspa ---
$ws.Range("E30").Value = "  -1.19%  "

# --- Row 31: InjectiveProtocol ---
Set-TextValue "D31" "34.76"
$ws.Range("E31").Value = "  -3.50%  "

# --- Row 32: OKB ---
Set-TextValue "D32" "49.58"
$ws.Range("E32").Value = "  -1.43%  "

# --- Row 33: Filecoin ---
Set-TextValue "D33" "5.50"
$ws.Range("E33").Value = "  +1.56%  "

# --- Row 34: Hedera ---
$ws.Range("E34").Value = "  -0.06%  "

# --- Row 35: FirstDigitalUSD ---
$ws.Range("E35").Value = "  -0.18%  "

# --- Row 36: Celestia ---
Set-TextValue "D36" "19.01"
$ws.Range("E36").Value = "  -2.88%  "

# --- Row 37: RenderToken ---
Set-TextValue "D37" "4.98"
$ws.Range("E37").Value = "  +4.66%  "

# --- Row 38: ARBITRUM ---
$ws.Range("E38").Value = "  -2.03%  "

# --- Row 39: LidoDAOToken ---
$ws.Range("E39").Value = "  +0.95%  "

# --- Row 40: Monero ---
Set-TextValue "D40" "128.88"
$ws.Range("E40").Value = "  +3.04%  "

# --- Row 41: EnergySwap ---
Set-TextValue "D41" "22.87"
$ws.Range("E41").Value = "  +2.99%  "

# --- Row 42: Stellar ---
$ws.Range("E42").Value = "  -1.20%  "

# --- Row 43: WEMIXToken ---
$ws.Range("E43").Value = "  -0.04%  "

# --- Row 44: VeChain ---
$ws.Range("E44").Value = "  +4.32%  "

# --- Row 45: Maker ---
Set-TextValue "D45" "2.053.67"
$ws.Range("E45").Value = "  -1.00%  "

# --- Rows 46/47: Stacks & NEARProtocol swap places (NEARProtocol now ranked 46th) ---
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D46" "3.26"
$ws.Range("E46").Value = "  -2.34%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D47" "2.15"
$ws.Range("E47").Value = "  +8.18%  "

# --- Row 48: ApeXProtocol ---
$ws.Range("E48").Value = "  -4.91%  "

# --- Row 49: FraxShare ---
Set-TextValue "D49" "8.86"
$ws.Range("E49").Value = "  -3.02%  "

# --- Row 50: THORChain ---
Set-TextValue "D50" "5.21"
$ws.Range("E50").Value = "  -3.42%  "

# --- Row 51: MultiversX ---
Set-TextValue "D51" "58.69"
$ws.Range("E51").Value = "  +0.12%  "
